# Natmi following Dr Hou advice
# Rebuild the LR-pairs table (Reln -> Vldlr) on Sheet1 with the updated
# sending/target cluster combinations (ECs, FAPs, sCs, M1, M2) and recomputed
# expression / specificity statistics produced by the new NATMI run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Reln"
$ws.Range("C2").Value = "Vldlr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.09865600000000001
$ws.Range("H2").Value = 0.295968
$ws.Range("I2").Value = 0.01009304870291239
$ws.Range("J2").Value = 0.01488643315542961
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.196272
$ws.Range("N2").Value = 0.588816
$ws.Range("O2").Value = 0.02813229386822481
$ws.Range("P2").Value = 0.03089009369338271
$ws.Range("Q2").Value = 0.019363410432
$ws.Range("R2").Value = 0.174270693888
$ws.Range("S2").Value = 0.0002839406121366367
$ws.Range("T2").Value = 0.0004598433149314993

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Reln"
$ws.Range("C3").Value = "Vldlr"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.09865600000000001
$ws.Range("H3").Value = 0.295968
$ws.Range("I3").Value = 0.01009304870291239
$ws.Range("J3").Value = 0.01488643315542961
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.808482333333333
$ws.Range("N3").Value = 14.425447
$ws.Range("O3").Value = 0.6892151609068061
$ws.Range("P3").Value = 0.7567787040415452
$ws.Range("Q3").Value = 0.4743856330773333
$ws.Range("R3").Value = 4.269470697696
$ws.Range("S3").Value = 0.006956282185817996
$ws.Range("T3").Value = 0.01126573559116711

# Row 4: ECs -> M1
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Reln"
$ws.Range("C4").Value = "Vldlr"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.09865600000000001
$ws.Range("H4").Value = 0.295968
$ws.Range("I4").Value = 0.01009304870291239
$ws.Range("J4").Value = 0.01488643315542961
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.03080033333333333
$ws.Range("N4").Value = 0.092401
$ws.Range("O4").Value = 0.004414710343669059
$ws.Range("P4").Value = 0.004847482995302872
$ws.Range("Q4").Value = 0.003038637685333333
$ws.Range("R4").Value = 0.027347739168
$ws.Range("S4").Value = 0.00004455788650790292
$ws.Range("T4").Value = 0.00007216173158165789

# Row 5: ECs -> M2
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Reln"
$ws.Range("C5").Value = "Vldlr"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.09865600000000001
$ws.Range("H5").Value = 0.295968
$ws.Range("I5").Value = 0.01009304870291239
$ws.Range("J5").Value = 0.01488643315542961
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.072589
$ws.Range("N5").Value = 0.217767
$ws.Range("O5").Value = 0.01040441366878908
$ws.Range("P5").Value = 0.01142435503336674
$ws.Range("Q5").Value = 0.007161340384000001
$ws.Range("R5").Value = 0.06445206345599999
$ws.Range("S5").Value = 0.0001050122538843356
$ws.Range("T5").Value = 0.0001700678975481098

# Row 6: ECs -> sCs
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Reln"
$ws.Range("C6").Value = "Vldlr"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.09865600000000001
$ws.Range("H6").Value = 0.295968
$ws.Range("I6").Value = 0.01009304870291239
$ws.Range("J6").Value = 0.01488643315542961
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.868607
$ws.Range("N6").Value = 3.737214
$ws.Range("O6").Value = 0.267833421212511
$ws.Range("P6").Value = 0.1960593642364025
$ws.Range("Q6").Value = 0.184349292192
$ws.Range("R6").Value = 1.106095753152
$ws.Range("S6").Value = 0.002703255764565523
$ws.Range("T6").Value = 0.002918624620201231

# Row 7: FAPs -> ECs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Reln"
$ws.Range("C7").Value = "Vldlr"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.2337746666666667
$ws.Range("H7").Value = 0.7013240000000001
$ws.Range("I7").Value = 0.02391642775070728
$ws.Range("J7").Value = 0.03527480283780177
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.196272
$ws.Range("N7").Value = 0.588816
$ws.Range("O7").Value = 0.02813229386822481
$ws.Range("P7").Value = 0.03089009369338271
$ws.Range("Q7").Value = 0.045883421376
$ws.Range("R7").Value = 0.4129507923840001
$ws.Range("S7").Value = 0.0006728239737610641
$ws.Range("T7").Value = 0.001089641964675299

# Row 8: FAPs -> FAPs
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Reln"
$ws.Range("C8").Value = "Vldlr"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.2337746666666667
$ws.Range("H8").Value = 0.7013240000000001
$ws.Range("I8").Value = 0.02391642775070728
$ws.Range("J8").Value = 0.03527480283780177
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 4.808482333333333
$ws.Range("N8").Value = 14.425447
$ws.Range("O8").Value = 0.6892151609068061
$ws.Range("P8").Value = 0.7567787040415452
$ws.Range("Q8").Value = 1.124101354647556
$ws.Range("R8").Value = 10.116912191828
$ws.Range("S8").Value = 0.01648356460051972
$ws.Range("T8").Value = 0.02669521957691264

# Row 9: FAPs -> M1
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Reln"
$ws.Range("C9").Value = "Vldlr"
$ws.Range("D9").Value = "M1"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.2337746666666667
$ws.Range("H9").Value = 0.7013240000000001
$ws.Range("I9").Value = 0.02391642775070728
$ws.Range("J9").Value = 0.03527480283780177
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.03080033333333333
$ws.Range("N9").Value = 0.092401
$ws.Range("O9").Value = 0.004414710343669059
$ws.Range("P9").Value = 0.004847482995302872
$ws.Range("Q9").Value = 0.007200337658222223
$ws.Range("R9").Value = 0.06480303892400001
$ws.Range("S9").Value = 0.0001055841009746612
$ws.Range("T9").Value = 0.0001709940069189055

# Row 10: FAPs -> M2
$ws.Range("A10").Value = "FAPs"
$ws.Range("B10").Value = "Reln"
$ws.Range("C10").Value = "Vldlr"
$ws.Range("D10").Value = "M2"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.2337746666666667
$ws.Range("H10").Value = 0.7013240000000001
$ws.Range("I10").Value = 0.02391642775070728
$ws.Range("J10").Value = 0.03527480283780177
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.072589
$ws.Range("N10").Value = 0.217767
$ws.Range("O10").Value = 0.01040441366878908
$ws.Range("P10").Value = 0.01142435503336674
$ws.Range("Q10").Value = 0.01696946927866667
$ws.Range("R10").Value = 0.152725223508
$ws.Range("S10").Value = 0.0002488364077980653
$ws.Range("T10").Value = 0.00040299187135106

# Row 11: FAPs -> sCs
$ws.Range("A11").Value = "FAPs"
$ws.Range("B11").Value = "Reln"
$ws.Range("C11").Value = "Vldlr"
$ws.Range("D11").Value = "sCs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.2337746666666667
$ws.Range("H11").Value = 0.7013240000000001
$ws.Range("I11").Value = 0.02391642775070728
$ws.Range("J11").Value = 0.03527480283780177
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 1.868607
$ws.Range("N11").Value = 3.737214
$ws.Range("O11").Value = 0.267833421212511
$ws.Range("P11").Value = 0.1960593642364025
$ws.Range("Q11").Value = 0.436832978556
$ws.Range("R11").Value = 2.620997871336
$ws.Range("S11").Value = 0.006405618667653768
$ws.Range("T11").Value = 0.00691595541794386

# Row 12: sCs -> ECs
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Reln"
$ws.Range("C12").Value = "Vldlr"
$ws.Range("D12").Value = "ECs"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 9.4422175
$ws.Range("H12").Value = 18.884435
$ws.Range("I12").Value = 0.9659905235463803
$ws.Range("J12").Value = 0.9498387640067686
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.196272
$ws.Range("N12").Value = 0.588816
$ws.Range("O12").Value = 0.02813229386822481
$ws.Range("P12").Value = 0.03089009369338271
$ws.Range("Q12").Value = 1.85324291316
$ws.Range("R12").Value = 11.11945747896
$ws.Range("S12").Value = 0.02717552928232711
$ws.Range("T12").Value = 0.02934060841377591

# Row 13: sCs -> FAPs
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Reln"
$ws.Range("C13").Value = "Vldlr"
$ws.Range("D13").Value = "FAPs"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 9.4422175
$ws.Range("H13").Value = 18.884435
$ws.Range("I13").Value = 0.9659905235463803
$ws.Range("J13").Value = 0.9498387640067686
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 4.808482333333333
$ws.Range("N13").Value = 14.425447
$ws.Range("O13").Value = 0.6892151609068061
$ws.Range("P13").Value = 0.7567787040415452
$ws.Range("Q13").Value = 45.40273603624083
$ws.Range("R13").Value = 272.416416217445
$ws.Range("S13").Value = 0.6657753141204683
$ws.Range("T13").Value = 0.7188177488734655

# Row 14: sCs -> M1
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Reln"
$ws.Range("C14").Value = "Vldlr"
$ws.Range("D14").Value = "M1"
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 9.4422175
$ws.Range("H14").Value = 18.884435
$ws.Range("I14").Value = 0.9659905235463803
$ws.Range("J14").Value = 0.9498387640067686
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.03080033333333333
$ws.Range("N14").Value = 0.092401
$ws.Range("O14").Value = 0.004414710343669059
$ws.Range("P14").Value = 0.004847482995302872
$ws.Range("Q14").Value = 0.2908234464058333
$ws.Range("R14").Value = 1.744940678435
$ws.Range("S14").Value = 0.004264568356186495
$ws.Range("T14").Value = 0.004604327256802308

# Row 15: sCs -> M2
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Reln"
$ws.Range("C15").Value = "Vldlr"
$ws.Range("D15").Value = "M2"
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 9.4422175
$ws.Range("H15").Value = 18.884435
$ws.Range("I15").Value = 0.9659905235463803
$ws.Range("J15").Value = 0.9498387640067686
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 0.3333333333333333
$ws.Range("M15").Value = 0.072589
$ws.Range("N15").Value = 0.217767
$ws.Range("O15").Value = 0.01040441366878908
$ws.Range("P15").Value = 0.01142435503336674
$ws.Range("Q15").Value = 0.6854011261075
$ws.Range("R15").Value = 4.112406756645
$ws.Range("S15").Value = 0.01005056500710668
$ws.Range("T15").Value = 0.01085129526446757

# Row 16: sCs -> sCs
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Reln"
$ws.Range("C16").Value = "Vldlr"
$ws.Range("D16").Value = "sCs"
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 9.4422175
$ws.Range("H16").Value = 18.884435
$ws.Range("I16").Value = 0.9659905235463803
$ws.Range("J16").Value = 0.9498387640067686
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 1.868607
$ws.Range("N16").Value = 3.737214
$ws.Range("O16").Value = 0.267833421212511
$ws.Range("P16").Value = 0.1960593642364025
$ws.Range("Q16").Value = 17.6437937160225
$ws.Range("R16").Value = 70.57517486409
$ws.Range("S16").Value = 0.2587245467802917
$ws.Range("T16").Value = 0.1862247841982574
